# Apply updated vm_pu values to the worksheet (case with 380 kV)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellValues = @{
    "B2" = 1.02
    "C2" = 1.018777190134748
    "D2" = 1.026513451475777
    "E2" = 1.028441947054263
    "F2" = 1.035878506188936
    "I2" = 1.028140619206577
    "J2" = 1.023983701230587
    "K2" = 1.029335878222202
    "L2" = 1.031258752783233
    "M2" = 1.038673849366234
    "N2" = 1.012006349856506
    "B3" = 1.02
    "C3" = 1.019837964448458
    "D3" = 1.026832767717716
    "E3" = 1.029421715039732
    "F3" = 1.037028821821524
    "I3" = 1.028119731977348
    "J3" = 1.024680361429499
    "K3" = 1.029464047931005
    "L3" = 1.032045993535909
    "M3" = 1.039632745861786
    "N3" = 1.012244501852445
    "B4" = 1.02
    "C4" = 1.020524907074266
    "D4" = 1.027038990663314
    "E4" = 1.030056559725998
    "F4" = 1.037774361322411
    "I4" = 1.02810454699596
    "J4" = 1.025131167901118
    "K4" = 1.029545830079758
    "L4" = 1.032555648281803
    "M4" = 1.040253847108459
    "N4" = 1.01239841392067
    "B5" = 1.02
    "C5" = 1.020813829861862
    "D5" = 1.027125589626558
    "E5" = 1.030323656259277
    "F5" = 1.038088075393573
    "I5" = 1.028097762610936
    "J5" = 1.025320691769214
    "K5" = 1.029579933844576
    "L5" = 1.032769968569184
    "M5" = 1.040515109076832
    "N5" = 1.012463073293099
    "B6" = 1.02
    "C6" = 1.020862349024004
    "D6" = 1.027140124207022
    "E6" = 1.030368515121524
    "F6" = 1.038140766328743
    "I6" = 1.028096599963333
    "J6" = 1.025352513939803
    "K6" = 1.029585643700015
    "L6" = 1.032805957485931
    "M6" = 1.040558984963684
    "N6" = 1.012473927226143
    "B7" = 1.02
    "C7" = 1.02052876715525
    "D7" = 1.027040148187686
    "E7" = 1.030060127866673
    "F7" = 1.037778552048582
    "I7" = 1.028104457918005
    "J7" = 1.025133700309516
    "K7" = 1.029546286868166
    "L7" = 1.032558511798797
    "M7" = 1.040257337510414
    "N7" = 1.012399278080471
    "B8" = 1.02
    "C8" = 1.019135569909541
    "D8" = 1.026621446480098
    "E8" = 1.028772884070614
    "F8" = 1.036267010527895
    "I8" = 1.02813390540539
    "J8" = 1.024219136428681
    "K8" = 1.029379431388346
    "L8" = 1.031524751091472
    "M8" = 1.038997782203664
    "N8" = 1.012086873282085
    "B9" = 1.02
    "C9" = 1.01668478420329
    "D9" = 1.025880711243438
    "E9" = 1.026511265982619
    "F9" = 1.033612736623003
    "I9" = 1.028173049635261
    "J9" = 1.022607722786081
    "K9" = 1.029076656696341
    "L9" = 1.029705111133826
    "M9" = 1.036783121677486
    "N9" = 1.011534943782861
    "B10" = 1.02
    "C10" = 1.015053735864303
    "D10" = 1.025385058414673
    "E10" = 1.025008019577426
    "F10" = 1.03184944187669
    "I10" = 1.028190631074268
    "J10" = 1.021533564086832
    "K10" = 1.02886901662911
    "L10" = 1.028493356457054
    "M10" = 1.035309935747672
    "N10" = 1.011166039416976
    "B11" = 1.02
    "C11" = 1.014348135187405
    "D11" = 1.025170031100778
    "E11" = 1.024358166769749
    "F11" = 1.03108739017691
    "I11" = 1.028196236062541
    "J11" = 1.021068470029094
    "K11" = 1.028777753015332
    "L11" = 1.027968971972987
    "M11" = 1.034672801649342
    "N11" = 1.011006076378505
    "B12" = 1.02
    "C12" = 1.014086141422113
    "D12" = 1.02509010175386
    "E12" = 1.024116942290496
    "F12" = 1.0308045503853
    "I12" = 1.028198017265132
    "J12" = 1.020895716942084
    "K12" = 1.028743651841878
    "L12" = 1.027774239312679
    "M12" = 1.03443625652106
    "N12" = 1.010946625336262
    "B13" = 1.02
    "C13" = 1.01414233555925
    "D13" = 1.025107249484421
    "E13" = 1.024168678527527
    "F13" = 1.030865210521971
    "I13" = 1.028197648785465
    "J13" = 1.020932772905002
    "K13" = 1.028750975765278
    "L13" = 1.027816007998315
    "M13" = 1.034486991058778
    "N13" = 1.010959379306702
    "B14" = 1.02
    "C14" = 1.014326476710061
    "D14" = 1.025163425298077
    "E14" = 1.024338223818738
    "F14" = 1.03106400606181
    "I14" = 1.028196389426981
    "J14" = 1.021054190132586
    "K14" = 1.028774938309027
    "L14" = 1.027952874349895
    "M14" = 1.034653246413062
    "N14" = 1.011001162822286
    "B15" = 1.02
    "C15" = 1.014439945087112
    "D15" = 1.025198029374264
    "E15" = 1.024442707445325
    "F15" = 1.03118651971999
    "I15" = 1.028195573671256
    "J15" = 1.021128999764434
    "K15" = 1.028789675730796
    "L15" = 1.028037208466125
    "M15" = 1.034755697038808
    "N15" = 1.011026902571408
    "B16" = 1.02
    "C16" = 1.015100577992576
    "D16" = 1.025399320714895
    "E16" = 1.025051170569725
    "F16" = 1.031900047620276
    "I16" = 1.028190216864643
    "J16" = 1.02156443134678
    "K16" = 1.028875045105079
    "L16" = 1.028528164730228
    "M16" = 1.035352236350497
    "N16" = 1.011176650909739
    "B17" = 1.02
    "C17" = 1.01551515029314
    "D17" = 1.025525478374927
    "E17" = 1.025433128034535
    "F17" = 1.032348017322074
    "I17" = 1.028186319765531
    "J17" = 1.021837572421634
    "K17" = 1.028928233685423
    "L17" = 1.02883621251074
    "M17" = 1.03572663446968
    "N17" = 1.01127052395519
    "B18" = 1.02
    "C18" = 1.015757026439381
    "D18" = 1.025599024633332
    "E18" = 1.025656019899017
    "F18" = 1.032609451984648
    "I18" = 1.02818385285959
    "J18" = 1.021996893281358
    "K18" = 1.02895912692504
    "L18" = 1.029015921795615
    "M18" = 1.035945088389231
    "N18" = 1.011325256793188
    "B19" = 1.02
    "C19" = 1.0158395106796
    "D19" = 1.025624095254975
    "E19" = 1.025732037653504
    "F19" = 1.03269861845784
    "I19" = 1.028182978809425
    "J19" = 1.02205121796653
    "K19" = 1.028969638492465
    "L19" = 1.029077203144091
    "M19" = 1.036019588056826
    "N19" = 1.011343915587419
    "B20" = 1.02
    "C20" = 1.015470664081779
    "D20" = 1.025511946906694
    "E20" = 1.025392137006533
    "F20" = 1.032299939776394
    "I20" = 1.028186757925733
    "J20" = 1.021808266732241
    "K20" = 1.028922540567381
    "L20" = 1.028803158772577
    "M20" = 1.035686457458765
    "N20" = 1.011260454511378
    "B21" = 1.02
    "C21" = 1.014272249036662
    "D21" = 1.025146884520153
    "E21" = 1.024288292544409
    "F21" = 1.031005459654759
    "I21" = 1.028196768571165
    "J21" = 1.021018435687276
    "K21" = 1.028767887494854
    "L21" = 1.027912569320712
    "M21" = 1.034604285185606
    "N21" = 1.010988859545931
    "B22" = 1.02
    "C22" = 1.013519323194119
    "D22" = 1.024917017289605
    "E22" = 1.023595186264539
    "F22" = 1.030192842195708
    "I22" = 1.028201323127522
    "J22" = 1.020521858121802
    "K22" = 1.028669484195791
    "L22" = 1.027352892817892
    "M22" = 1.033924544930364
    "N22" = 1.010817902416904
    "B23" = 1.02
    "C23" = 1.013918410019024
    "D23" = 1.025038905493229
    "E23" = 1.023962527342228
    "F23" = 1.030623505429295
    "I23" = 1.028199073251911
    "J23" = 1.020785101339317
    "K23" = 1.028721759700743
    "L23" = 1.027649562091004
    "M23" = 1.034284825157511
    "N23" = 1.010908548417934
    "B24" = 1.02
    "C24" = 1.015490765290502
    "D24" = 1.025518061316702
    "E24" = 1.025410658770727
    "F24" = 1.03232166351239
    "I24" = 1.028186560538906
    "J24" = 1.021821508704644
    "K24" = 1.028925113446604
    "L24" = 1.028818094240572
    "M24" = 1.035704611492333
    "N24" = 1.011265004526705
    "B25" = 1.02
    "C25" = 1.017317874621407
    "D25" = 1.026072541693807
    "E25" = 1.027095156568681
    "F25" = 1.034297833428938
    "I25" = 1.028164434695646
    "J25" = 1.023024291679433
    "K25" = 1.029155958733625
    "L25" = 1.030175296368112
    "M25" = 1.037355091090547
    "N25" = 1.011677799327959
}

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
